# The commit swaps the two theme parts of the deck:
#   ppt/theme/theme1.xml  (the slide master's theme, previously the
#                          "Integral" colour scheme)
#   ppt/theme/theme2.xml  (the notes master's theme, previously the
#                          plain "Office Theme" colour scheme)
# become each other - theme1.xml ends up holding the "Office Theme"
# colours, theme2.xml ends up holding the "Integral" colours. The font
# scheme and format scheme (gradients/line styles/effects) are already
# byte-identical between the two theme parts, so the only real content
# difference is the <a:clrScheme> (and the cosmetic name attributes).
#
# The PowerPoint object model only exposes the slide-master's theme
# colours for writing (SlideMaster.Theme.ThemeColorScheme); the notes
# master's theme colours are not independently writable through COM, so
# this script applies the reachable half of the swap: it recolors the
# master theme (theme1.xml) to the "Office Theme" palette.

$p = $ppt.ActivePresentation
$scheme = $p.SlideMaster.Theme.ThemeColorScheme

# Target palette ("Office Theme"), in ThemeColorScheme.Item() order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink.
# RGB() on this host is the usual COLORREF ordering (B*65536 + G*256 + R),
# so each entry below is built from the R,G,B bytes of the target hex
# colour rather than used as a plain 0xRRGGBB literal.
$targets = @(
    @{ Index = 1;  R = 0x00; G = 0x00; B = 0x00 },  # dk1      000000
    @{ Index = 2;  R = 0xFF; G = 0xFF; B = 0xFF },  # lt1      FFFFFF
    @{ Index = 3;  R = 0x44; G = 0x54; B = 0x6A },  # dk2      44546A
    @{ Index = 4;  R = 0xE7; G = 0xE6; B = 0xE6 },  # lt2      E7E6E6
    @{ Index = 5;  R = 0x5B; G = 0x9B; B = 0xD5 },  # accent1  5B9BD5
    @{ Index = 6;  R = 0xED; G = 0x7D; B = 0x31 },  # accent2  ED7D31
    @{ Index = 7;  R = 0xA5; G = 0xA5; B = 0xA5 },  # accent3  A5A5A5
    @{ Index = 8;  R = 0xFF; G = 0xC0; B = 0x00 },  # accent4  FFC000
    @{ Index = 9;  R = 0x44; G = 0x72; B = 0xC4 },  # accent5  4472C4
    @{ Index = 10; R = 0x70; G = 0xAD; B = 0x47 },  # accent6  70AD47
    @{ Index = 11; R = 0x05; G = 0x63; B = 0xC1 },  # hlink    0563C1
    @{ Index = 12; R = 0x95; G = 0x4F; B = 0x72 }   # folHlink 954F72
)

foreach ($t in $targets) {
    $colorRef = ($t.B * 65536) + ($t.G * 256) + $t.R
    $scheme.Item($t.Index).RGB = $colorRef
}
